$wb = $excel.ActiveWorkbook

# --- DLC_List sheet: pregenerated-BTS wind cases ---
$dlc = $wb.Worksheets.Item("DLC_List")
$dlc.Activate()

# Row 3 (1p3 / ETM): narrow wind-speed bin and the BlPitch/Yaw pregenerated-BTS ranges
$dlc.Range("C3").Value = "[10 12]"
$dlc.Range("L3").Value = "<b>[5 7]"
$dlc.Range("M3").Value = "<b>[5 7]"
$dlc.Range("N3").Value = "<a>[-2:2:2]"
$dlc.Range("O3").Value = "<a>[-2:2:2]"
$dlc.Range("P3").Value = "<a>[-2:2:2]"

# New row 8: spatially-coherent pregenerated BTS wind case
$dlc.Range("A8").Value = "coh"
$dlc.Range("B8").Value = "BTS:NTM_URef-%d_turbsim_coh"
$dlc.Range("C8").Value = "[10 20]"

# Move the active selection to the newly edited cell
[void]$dlc.Range("N3").Select()

# --- config sheet: selection unchanged (still B14) ---
$cfg = $wb.Worksheets.Item("config")
[void]$cfg.Activate()
[void]$cfg.Range("B14").Select()

[void]$dlc.Activate()
